$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 15:22"

# Serbia (row 42) - refreshed case counts
$ws.Range("B42").Value = 8497
$ws.Range("C42").Value = 222
$ws.Range("D42").Value = 1260
$ws.Range("E42").Value = 7069
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 168

# Rows 73/74: Camerun and Azerbaiyan swap places (alphabetical re-sort) and
# get refreshed figures for the row that is now Azerbaiyan.
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 1717
$ws.Range("C73").Value = 39
$ws.Range("D73").Value = 1221
$ws.Range("E73").Value = 474
$ws.Range("F73").Value = 15
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 22

$ws.Range("A74").Value = "Camerun"
$ws.Range("B74").Value = 1705
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 805
$ws.Range("E74").Value = 842
$ws.Range("F74").Value = 12
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 58

# Rows 79/80: Eslovenia and Republica de Macedonia swap places and get
# refreshed figures for the row that is now Republica de Macedonia.
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 1421
$ws.Range("C79").Value = 22
$ws.Range("D79").Value = 589
$ws.Range("E79").Value = 761
$ws.Range("F79").Value = 13
$ws.Range("G79").Value = 6
$ws.Range("H79").Value = 71

$ws.Range("A80").Value = "Eslovenia"
$ws.Range("B80").Value = 1408
$ws.Range("C80").Value = 6
$ws.Range("D80").Value = 223
$ws.Range("E80").Value = 1099
$ws.Range("F80").Value = 24
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 86
